$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '26.915.28'
$ws.Cells.Item(2, 5).Value = '  +0.08%  '
$ws.Cells.Item(3, 4).Value = '1.549.77'
$ws.Cells.Item(3, 5).Value = '  -0.13%  '
$ws.Cells.Item(4, 5).Value = '  -0.38%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '206.42'
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +0.04%  '
$ws.Cells.Item(6, 5).Value = '  +1.01%  '
$ws.Cells.Item(7, 5).Value = '  -0.37%  '
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '22.07'
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  +2.63%  '
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.246'
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  -0.20%  '
$ws.Cells.Item(10, 5).Value = '  +0.86%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.0856'
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  -0.14%  '
$ws.Cells.Item(12, 4).Value = '1.771.61'
$ws.Cells.Item(13, 4).Value = '1.550.43'
$ws.Cells.Item(13, 5).Value = '  -2.23%  '
$ws.Cells.Item(14, 5).Value = '  +0.92%  '
$ws.Cells.Item(15, 5).Value = '  +0.93%  '
$ws.Cells.Item(16, 4).Value = '26.898.06'
$ws.Cells.Item(16, 5).Value = '  -0.02%  '
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '61.65'
$ws.Cells.Item(17, 4).Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  +0.06%  '
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '217.01'
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  +1.55%  '
$ws.Cells.Item(19, 4).Value = '0.0₃0700'
$ws.Cells.Item(19, 5).Value = '  +2.20%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '7.26'
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  +0.43%  '
$ws.Cells.Item(21, 5).Value = '  -0.39%  '
$ws.Cells.Item(22, 5).Value = '  +0.32%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '9.21'
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  +0.61%  '
$ws.Cells.Item(24, 5).Value = '  -0.99%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '153.86'
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  +0.58%  '
$ws.Cells.Item(26, 5).Value = '  -0.41%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '14.95'
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  +0.66%  '
$ws.Cells.Item(28, 5).Value = '  +0.82%  '
$ws.Cells.Item(29, 5).Value = '  -0.34%  '
$ws.Cells.Item(30, 5).Value = '  +1.97%  '
$ws.Cells.Item(31, 5).Value = '  -0.52%  '
$ws.Cells.Item(32, 5).Value = '  -0.16%  '
$ws.Cells.Item(33, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '3.10'
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  +4.78%  '
$ws.Cells.Item(34, 2).Value = 'Maker'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(34, 4).Value = '1.413.93'
$ws.Cells.Item(34, 5).Value = '  +3.29%  '
$ws.Cells.Item(35, 5).Value = '  +2.68%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.968'
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  -0.54%  '
$ws.Cells.Item(37, 5).Value = '  +0.12%  '
$ws.Cells.Item(38, 5).Value = '  +1.09%  '
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.529'
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  +1.10%  '
$ws.Cells.Item(40, 5).Value = '  +0.01%  '
$ws.Cells.Item(41, 5).Value = '  -0.32%  '
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '5.69'
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  +3.31%  '
$ws.Cells.Item(43, 5).Value = '  +2.92%  '
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '1.00'
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  +1.48%  '
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '64.63'
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  +1.75%  '
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '1.74'
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  +0.50%  '
$ws.Cells.Item(47, 4).Value = '1.685.16'
$ws.Cells.Item(47, 5).Value = '  -0.08%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '87.46'
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  +1.48%  '
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '0.0516'
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  +1.91%  '
$ws.Cells.Item(50, 5).Value = '  +3.57%  '
$ws.Cells.Item(51, 5).Value = '  +0.76%  '
